$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'23.442.50"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "
$c = $ws.Range("D3")
$c.Value = "'1.632.94"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  +0.37%  "
$c = $ws.Range("D5")
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$c = $ws.Range("D6")
$c.Value = "'304.83"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.27%  "
$c = $ws.Range("D7")
$c.Value = "'0.3759"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$c = $ws.Range("D8")
$c.Value = "'0.3637"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D9")
$c.Value = "'51.62"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.49%  "
$c = $ws.Range("D10")
$c.Value = "'0.08181"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.29%  "
$c = $ws.Range("D11")
$c.Value = "'1.223"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.64%  "
$c = $ws.Range("D12")
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
$c = $ws.Range("D13")
$c.Value = "'22.56"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.56%  "
$c = $ws.Range("D14")
$c.Value = "'6.549"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.04%  "
$c = $ws.Range("D15")
$c.Value = "'0.00001250"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.67%  "
$c = $ws.Range("D16")
$c.Value = "'7.247"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.23%  "
$c = $ws.Range("D17")
$c.Value = "'1.630.22"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.82%  "
$c = $ws.Range("D18")
$c.Value = "'94.45"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.99%  "
$c = $ws.Range("D19")
$c.Value = "'0.06980"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.50%  "
$c = $ws.Range("D20")
$c.Value = "'17.77"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.21%  "
$c = $ws.Range("D21")
$c.Value = "'6.459"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.16%  "
$c = $ws.Range("D22")
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.27%  "
$c = $ws.Range("D23")
$c.Value = "'12.74"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.31%  "
$c = $ws.Range("D24")
$c.Value = "'23.452.15"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.50%  "
$c = $ws.Range("D25")
$c.Value = "'3.182"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.84%  "
$c = $ws.Range("D26")
$c.Value = "'2.473"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.24%  "
$c = $ws.Range("D27")
$c.Value = "'21.32"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.23%  "
$c = $ws.Range("D28")
$c.Value = "'149.85"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.93%  "
$c = $ws.Range("D29")
$c.Value = "'5.312"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.47%  "
$c = $ws.Range("D30")
$c.Value = "'134.90"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.23%  "
$c = $ws.Range("D31")
$c.Value = "'1.811.98"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.60%  "
$c = $ws.Range("D32")
$c.Value = "'2.287"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.84%  "
$c = $ws.Range("D33")
$c.Value = "'6.820"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.81%  "
$c = $ws.Range("D34")
$c.Value = "'1.017"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.59%  "
$c = $ws.Range("D35")
$c.Value = "'11.21"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +6.95%  "
$c = $ws.Range("D36")
$c.Value = "'0.02788"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.13%  "
$c = $ws.Range("D37")
$c.Value = "'0.2529"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D38")
$c.Value = "'0.08795"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D39")
$c.Value = "'0.07165"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.69%  "
$c = $ws.Range("D40")
$c.Value = "'6.065"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.47%  "
$c = $ws.Range("D41")
$c.Value = "'0.7069"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.19%  "
$c = $ws.Range("D42")
$c.Value = "'1.348"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.65%  "
$c = $ws.Range("D43")
$c.Value = "'16.19"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.19%  "
$c = $ws.Range("D44")
$c.Value = "'12.29"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.77%  "
$c = $ws.Range("D45")
$c.Value = "'0.6516"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.20%  "
$c = $ws.Range("D46")
$c.Value = "'2.330"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.32%  "
$c = $ws.Range("D47")
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.16%  "
$c = $ws.Range("D48")
$c.Value = "'3.992"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.18%  "
$c = $ws.Range("D49")
$c.Value = "'0.08041"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.10%  "
$c = $ws.Range("D50")
$c.Value = "'1.207"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.79%  "
$c = $ws.Range("D51")
$c.Value = "'125.06"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.72%  "
